$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for the extra columns
$ws.Range("C1").Value = "WIN"
$ws.Range("D1").Value = "TOP4"
$ws.Range("E1").Value = "TOP5"
$ws.Range("F1").Value = "TOP6"
$ws.Range("G1").Value = "RELEGATION"
$ws.Range("H1").Value = "ExpPoints"

# Copy header style (bold/border/centered) from B1 to the new header cells
$ws.Range("B1").Copy()
$ws.Range("C1:H1").PasteSpecial(-4122)  # xlPasteFormats

# New team ordering + ExpPoints (now in column H), with blank placeholder cells in C:G
$teams = @(
    @{ Row = 2;  Team = "Barcelona";           Exp = 86.8386200201538 },
    @{ Row = 3;  Team = "Real Madrid";         Exp = 86.04531533110264 },
    @{ Row = 4;  Team = "Atlético de Madrid";  Exp = 70.00314956420127 },
    @{ Row = 5;  Team = "Villarreal";          Exp = 63.86196922548655 },
    @{ Row = 6;  Team = "Real Betis";          Exp = 62.48926295936266 },
    @{ Row = 7;  Team = "Athletic Club";       Exp = 56.68612634844217 },
    @{ Row = 8;  Team = "Rayo Vallecano";      Exp = 52.9486176969384 },
    @{ Row = 9;  Team = "Sevilla";             Exp = 49.21931099513707 },
    @{ Row = 10; Team = "Getafe";              Exp = 48.79097228784917 },
    @{ Row = 11; Team = "Celta de Vigo";       Exp = 48.19862895940904 },
    @{ Row = 12; Team = "Espanyol";            Exp = 47.74023850437959 },
    @{ Row = 13; Team = "Osasuna";             Exp = 46.73885588594264 },
    @{ Row = 14; Team = "Valencia";            Exp = 45.95150206962512 },
    @{ Row = 15; Team = "Real Sociedad";       Exp = 45.01540212389604 },
    @{ Row = 16; Team = "Elche";               Exp = 40.69738607757875 },
    @{ Row = 17; Team = "Alavés";              Exp = 39.96524084463713 },
    @{ Row = 18; Team = "Mallorca";            Exp = 39.07289565114249 },
    @{ Row = 19; Team = "Levante";             Exp = 34.76243244862623 },
    @{ Row = 20; Team = "Real Oviedo";         Exp = 32.16607011582397 },
    @{ Row = 21; Team = "Girona";              Exp = 31.59460256962696 }
)

foreach ($t in $teams) {
    $r = $t.Row
    $ws.Range("B$r").Value = $t.Team
    $ws.Range("C$r").Value = ""
    $ws.Range("D$r").Value = ""
    $ws.Range("E$r").Value = ""
    $ws.Range("F$r").Value = ""
    $ws.Range("G$r").Value = ""
    $ws.Range("H$r").Value = $t.Exp
}
